# Bump cell numbers in the operator table to match 1-based numbering used in the adoc.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("F2").Value2 = 1
$ws.Range("G2").Value2 = 2
$ws.Range("H2").Value2 = 3
$ws.Range("I2").Value2 = 4
$ws.Range("J2").Value2 = 5
$ws.Range("K2").Value2 = 6
$ws.Range("L2").Value2 = 7
$ws.Range("M2").Value2 = 8

# Row 3
$ws.Range("F3").Value2 = 9
$ws.Range("G3").Value2 = 10
$ws.Range("H3").Value2 = 11
$ws.Range("I3").Value2 = 12
$ws.Range("J3").Value2 = 13
$ws.Range("K3").Value2 = 14
$ws.Range("L3").Value2 = 15
$ws.Range("M3").Value2 = 16

# Row 4
$ws.Range("F4").Value2 = 17
$ws.Range("G4").Value2 = 18
$ws.Range("H4").Value2 = 19
$ws.Range("I4").Value2 = 20
$ws.Range("J4").Value2 = 21
$ws.Range("K4").Value2 = 22
$ws.Range("L4").Value2 = 23
$ws.Range("M4").Value2 = 24

# Row 5
$ws.Range("F5").Value2 = 25
$ws.Range("G5").Value2 = 26
$ws.Range("H5").Value2 = 27
$ws.Range("I5").Value2 = 28
$ws.Range("J5").Value2 = 29
$ws.Range("K5").Value2 = 30
$ws.Range("L5").Value2 = 31
$ws.Range("M5").Value2 = 32

# Row 6
$ws.Range("F6").Value2 = 33
$ws.Range("G6").Value2 = 34
$ws.Range("H6").Value2 = 35
$ws.Range("I6").Value2 = 36
$ws.Range("J6").Value2 = 37
$ws.Range("K6").Value2 = 38
$ws.Range("L6").Value2 = 39
$ws.Range("M6").Value2 = 40

# Row 7
$ws.Range("F7").Value2 = 41
$ws.Range("G7").Value2 = 42
$ws.Range("H7").Value2 = 43
$ws.Range("I7").Value2 = 44
$ws.Range("J7").Value2 = 45
$ws.Range("K7").Value2 = 46
$ws.Range("L7").Value2 = 47
$ws.Range("M7").Value2 = 48

# Row 8
$ws.Range("F8").Value2 = 49
$ws.Range("G8").Value2 = 50
$ws.Range("H8").Value2 = 51
$ws.Range("I8").Value2 = 52
$ws.Range("J8").Value2 = 53
$ws.Range("K8").Value2 = 54
$ws.Range("L8").Value2 = 55
$ws.Range("M8").Value2 = 56

# Row 9
$ws.Range("F9").Value2 = 57
$ws.Range("G9").Value2 = 58
$ws.Range("H9").Value2 = 59
$ws.Range("I9").Value2 = 60
$ws.Range("J9").Value2 = 61
$ws.Range("K9").Value2 = 62
$ws.Range("L9").Value2 = 63
$ws.Range("M9").Value2 = 64

# Row 15
$ws.Range("M15").Value2 = 8

# Row 16
$ws.Range("F16").Value2 = 9
$ws.Range("H16").Value2 = 11
$ws.Range("J16").Value2 = 13
$ws.Range("L16").Value2 = 15
$ws.Range("M16").Value2 = 16
$ws.Range("P16").Value2 = 9

# Row 17
$ws.Range("P17").Value2 = 32

# Row 19
$ws.Range("M19").Value2 = 40
$ws.Range("P19").Value2 = 11

# Row 20
$ws.Range("F20").Value2 = 41
$ws.Range("G20").Value2 = 42
$ws.Range("J20").Value2 = 45
$ws.Range("K20").Value2 = 46
$ws.Range("M20").Value2 = 48
$ws.Range("P20").Value2 = 15

# Row 21
$ws.Range("P21").Value2 = 58

# Row 22
$ws.Range("F22").Value2 = 57
$ws.Range("G22").Value2 = 58
$ws.Range("K22").Value2 = 62
$ws.Range("O22").Value2 = 16
$ws.Range("P22").Value2 = 57

# Row 23
$ws.Range("P23").Value2 = 41

# Row 24
$ws.Range("P24").Value2 = 13

# Row 25
$ws.Range("P25").Value2 = 40

# Row 26
$ws.Range("P26").Value2 = 48

# Row 27
$ws.Range("P27").Value2 = 62

# Row 28
$ws.Range("P28").Value2 = 45

# Update the saved cell selection to match the author's final cursor position.
$ws.Range("P29").Select()
